# Finish Toronto surprise songs — populate the remaining "song"/"mashup"/
# "guest"/"dress" columns for surprise-song rows 282-293 (Toronto, Leg 2),
# matching the source data update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 282 (guitar): dress + song + mashup ---
$ws.Range("E282").Value = "koi fish"
$ws.Range("G282").Value = "My Boy Only Breaks His Favorite Toys"
$ws.Range("H282").Value = "This Is Why We Can't Have Nice Things"

# --- Row 283 (piano): dress + song (quote-prefixed) + mashup ---
$ws.Range("E283").Value = "koi fish"
$ws.Range("G283").Value = "False God"
$ws.Range("H283").Value = "''tis the damn season"

# --- Row 284 (guitar): dress + song (quote-prefixed) + mashup ---
$ws.Range("E284").Value = "supernova"
$ws.Range("G284").Value = "'I Don't Wanna Live Forever"
$ws.Range("H284").Value = "Mine (Taylor's Version)"

# --- Row 285 (piano): dress + song + mashup ---
$ws.Range("E285").Value = "supernova"
$ws.Range("G285").Value = "evermore"
$ws.Range("H285").Value = "Peter"

# --- Row 286 (guitar): dress + song + mashup + guest ---
$ws.Range("E286").Value = "sunrise boulevard"
$ws.Range("G286").Value = "us."
$ws.Range("H286").Value = "Out Of The Woods (Taylor's Version)"
$ws.Range("I286").Value = "Gracie Abrams"

# --- Row 287 (piano): dress + song + mashup ---
$ws.Range("E287").Value = "sunrise boulevard"
$ws.Range("G287").Value = "You're On Your Own, Kid"
$ws.Range("H287").Value = "long story short"

# --- Row 288 (guitar): dress + song + mashup ---
$ws.Range("E288").Value = "betta fish"
$ws.Range("G288").Value = "Mr. Perfectly Fine (Taylor's Version) [From The Vault]"
$ws.Range("H288").Value = "Better Than Revenge (Taylor's Version)"

# --- Row 289 (piano): dress + song + mashup ---
$ws.Range("E289").Value = "betta fish"
$ws.Range("G289").Value = "State Of Grace (Taylor's Version)"
$ws.Range("H289").Value = "Labyrinth"

# --- Row 290 (guitar): dress + song + mashup ---
$ws.Range("E290").Value = "koi fish"
$ws.Range("G290").Value = "Ours (Taylor's Version)"
$ws.Range("H290").Value = "the last great american dynasty"

# --- Row 291 (piano): dress + song + mashup ---
$ws.Range("E291").Value = "koi fish"
$ws.Range("G291").Value = "Cassandra"
$ws.Range("H291").Value = "mad woman; I Did Something Bad"

# --- Row 292 (guitar): dress + song + mashup ---
$ws.Range("E292").Value = "supernova"
$ws.Range("G292").Value = "Sparks Fly (Taylor's Version)"
$ws.Range("H292").Value = "Message In A Bottle (Taylor's Version) [From The Vault]"

# --- Row 293 (piano): dress + song + mashup ---
$ws.Range("E293").Value = "supernova"
$ws.Range("G293").Value = "You're Losing Me (From The Vault)"
$ws.Range("H293").Value = "How Did It End?"

# Leave the view parked on the last entry, like the source edit did.
$ws.Range("G291").Select() | Out-Null
